$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts existing rows 13..197 down to 14..198),
# matching the weekly price-list pattern of a newly recorded observation.
$ws.Rows(13).Insert()

# Fill the new row 13 with the standard Femacal de La Calera / Ciboulette record,
# dated 2021-10-27 (serial 44496) with a volume of 150.
$ws.Range("A13").Value = 3
$ws.Range("B13").Value = "Femacal de La Calera"
$ws.Range("C13").Value = "Coquimbo"
$ws.Range("D13").Value = 44496
$ws.Range("D13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E13").Value = 5
$ws.Range("F13").Value = 100112039
$ws.Range("G13").Value = "Ciboulette"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 150
$ws.Range("K13").Value = 1500
$ws.Range("L13").Value = 1500
$ws.Range("M13").Value = 1500
$ws.Range("N13").Value = "$/docena de atados"
$ws.Range("O13").Value = "Provincia de Quillota"
$ws.Range("P13").Value = 500
$ws.Range("Q13").Value = 3
$ws.Range("R13").Value = "Hortaliza"
